# Apply the FrMethodOfAdministration ValueSet spreadsheet updates:
#  - Update the canonical URL to the new HL7 France location
#  - Update the build/generation Date to the new timestamp
#  - Clear the Copyright cell (value removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/ValueSet/FrMethodOfAdministration"
$ws.Range("B8").Value = "2025-04-10T15:35:36+00:00"
$ws.Range("B14").ClearContents()
